# aggiornamento fino al 26/03
# Append 5 new daily rows (234-238) to the data table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, date serial (A), nuovi pos. (B), somma mobile 7gg. (C), somma mobile 7gg. per 100mila abitanti (D)
$newRows = @(
    @(234, 44308, 9, 43, 125.1236687423616),
    @(235, 44309, 6, 42, 122.2138159809114),
    @(236, 44310, 7, 41, 119.3039632194611),
    @(237, 44311, 3, 37, 107.66455217366),
    @(238, 44312, 8, 39, 113.4842576965606)
)

# Last pre-existing row (233) carries the date-column style (bold, centered,
# bordered, custom date-time number format) that the new rows must reuse.
$lastRow = 233

foreach ($entry in $newRows) {
    $r = $entry[0]

    # Copy the formatting of column A from the previous row so the new date
    # cell keeps the same style (border/alignment/number format) as the rest
    # of the column.
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]

    $lastRow = $r
}
